{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that ends the \"\u0420\u0430\u0437\u043c\u0435\u0440 \u0438 \u043f\u043e\u043b\u043e\u0436\u0435\u043d\u0438\u0435 \u043d\u0430 \u043f\u0440\u043e\u0437\u043e\u0440\u0435\u0446\u0430\" example\n// (the one containing the window.setSize(...) code line) -- the new Selenium\n// tricks section is appended immediately after it, before the trailing\n// HTMLBottomofForm paragraph that closes the document body.\nconst anchorText = \"window.setSize(new Dimension(800, 600));\";\nlet anchorParagraph = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(anchorText) >= 0) {\n    anchorParagraph = p;\n  }\n}\nif (!anchorParagraph) {\n  throw new Error(\"Anchor paragraph not found\");\n}\n\n// Insert a fresh empty paragraph right after the anchor, then replace its\n// (collapsed) range with the OOXML for the full block of new paragraphs --\n// Word resolves the insertion point relative to that new paragraph, so the\n// content lands exactly between the anchor paragraph and HTMLBottomofForm.\nconst newParagraph = anchorParagraph.insertParagraph(\"\", Word.InsertLocation.after);\nawait context.sync();\n\nconst blockOoxml = '<w:p><w:pPr><w:pStyle w:val=\"Heading2\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:before=\"150\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0438\u0440</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0430\u043d\u0435 \u043d\u0430</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> \u0434\u0435\u0439\u0441\u0442\u0432\u0438</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u044f \u0432</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> Selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>EventFiringWebDriver: Java</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>public static class MyListener extends AbstractWebDriverEventListener {</w:t><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void beforeFindBy(By by, WebElement element, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(by);</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void afterFindBy(By by, WebElement element, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(by +  \" found\");</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void onException(Throwable throwable, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(throwable);</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>driver = new EventFiringWebDriver(new ChromeDriver());</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0421\u043a\u0440\u0438\u043d\u0448\u043e\u0442</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>File tempFile = ((TakesScreenshot) driver).getScreenshotAs(OutputType.FILE);</w:t><w:br/><w:t>try {</w:t><w:br/><w:t xml:space=\"preserve\">  Files.copy(tempFile, new File(\"screen.png\"));</w:t><w:br/><w:t>} catch (IOException e) {</w:t><w:br/><w:t xml:space=\"preserve\">  e.printStackTrace();</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0414\u043e\u0441\u0442\u044a\u043f \u043a\u044a\u043c \u043b\u043e\u0433\u043e\u0432\u0435\u0442\u0435 \u043d\u0430 \u0411\u0440\u0430\u0443\u0437\u044a\u0440\u0430</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>for (LogEntry l : driver.manage().logs().get(\"browser\").getAll()) {</w:t><w:br/><w:t xml:space=\"preserve\">    System.out.println(l);</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:before=\"150\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>';\nconst flatOpc = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + blockOoxml + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nconst targetRange = newParagraph.getRange(\"Whole\");\ntargetRange.insertOoxml(flatOpc, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph ending the \"window.setSize(...)\" code example -- the\n# new Selenium tricks section is appended right after it, before the\n# trailing HTMLBottomofForm paragraph that closes the document body.\n$anchorText = \"window.setSize(new Dimension(800, 600));\"\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*$anchorText*\") {\n        $anchor = $p\n    }\n}\nif ($anchor -eq $null) {\n    throw \"Anchor paragraph not found\"\n}\n\n# Insert a fresh empty paragraph right after the anchor, then replace its\n# content with the OOXML for the full block of new paragraphs.\n$anchor.Range.InsertParagraphAfter() | Out-Null\n$newPara = $anchor.Next()\n\n$blockXml = '<w:p><w:pPr><w:pStyle w:val=\"Heading2\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:before=\"150\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>\u041f\u0440\u043e\u0442\u043e\u043a\u043e\u043b\u0438\u0440</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0430\u043d\u0435 \u043d\u0430</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> \u0434\u0435\u0439\u0441\u0442\u0432\u0438</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u044f \u0432</w:t></w:r><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t xml:space=\"preserve\"> Selenium</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>EventFiringWebDriver: Java</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>public static class MyListener extends AbstractWebDriverEventListener {</w:t><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void beforeFindBy(By by, WebElement element, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(by);</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void afterFindBy(By by, WebElement element, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(by +  \" found\");</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:br/><w:t xml:space=\"preserve\">    @Override</w:t><w:br/><w:t xml:space=\"preserve\">    public void onException(Throwable throwable, WebDriver driver) {</w:t><w:br/><w:t xml:space=\"preserve\">        System.out.println(throwable);</w:t><w:br/><w:t xml:space=\"preserve\">    }</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>driver = new EventFiringWebDriver(new ChromeDriver());</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0421\u043a\u0440\u0438\u043d\u0448\u043e\u0442</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>File tempFile = ((TakesScreenshot) driver).getScreenshotAs(OutputType.FILE);</w:t><w:br/><w:t>try {</w:t><w:br/><w:t xml:space=\"preserve\">  Files.copy(tempFile, new File(\"screen.png\"));</w:t><w:br/><w:t>} catch (IOException e) {</w:t><w:br/><w:t xml:space=\"preserve\">  e.printStackTrace();</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading3\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:beforeAutospacing=\"0\" w:before=\"150\" w:afterAutospacing=\"0\" w:after=\"150\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b w:val=\"false\"/><w:bCs w:val=\"false\"/><w:color w:val=\"333333\"/><w:sz w:val=\"36\"/><w:szCs w:val=\"36\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u0414\u043e\u0441\u0442\u044a\u043f \u043a\u044a\u043c \u043b\u043e\u0433\u043e\u0432\u0435\u0442\u0435 \u043d\u0430 \u0411\u0440\u0430\u0443\u0437\u044a\u0440\u0430</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading4\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"300\" w:after=\"300\"/><w:rPr><w:lang w:val=\"bg-BG\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:cs=\"Arial\" w:ascii=\"Arial\" w:hAnsi=\"Arial\"/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"30\"/><w:szCs w:val=\"30\"/><w:lang w:val=\"bg-BG\"/></w:rPr><w:t>\u041f\u0440\u0438\u043c\u0435\u0440</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"HTMLPreformatted\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"F5F5F5\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"300\" w:before=\"0\" w:after=\"150\"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Consolas\" w:hAnsi=\"Consolas\"/><w:color w:val=\"333333\"/><w:lang w:val=\"en-US\"/></w:rPr><w:t>for (LogEntry l : driver.manage().logs().get(\"browser\").getAll()) {</w:t><w:br/><w:t xml:space=\"preserve\">    System.out.println(l);</w:t><w:br/><w:t>}</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"Heading2\"/><w:shd w:val=\"clear\" w:color=\"auto\" w:fill=\"FFFFFF\"/><w:spacing w:lineRule=\"atLeast\" w:line=\"600\" w:before=\"150\" w:after=\"150\"/><w:rPr><w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:cs=\"Arial\"/><w:b/><w:b/><w:bCs/><w:color w:val=\"333333\"/><w:sz w:val=\"42\"/><w:szCs w:val=\"42\"/><w:lang w:val=\"en-US\"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p>'\n$flatOpc = '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + $blockXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$newPara.Range.InsertXML($flatOpc) | Out-Null\n"}
